# Commit: "Get a valid Constellation name from the User to show the date
# asked by the user." (see #2)
#
# The "North_cons" table listed the constellation "Leo" in lower/title case;
# it needs to be the valid, all-caps constellation abbreviation "LEO" so
# look-ups against the user-entered name succeed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("North")

# Fix the constellation name in the North table (row 5 -> "Leo" becomes "LEO").
$ws.Range("A5").Value = "LEO"

# Leave the selection where the editor ended up after making the change.
$ws.Range("A14").Select()
